$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (same bold/border/centered style used by
# every other header cell) onto the new H1 header, then set its text.
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 8).Value = "Save"

# New "Save" column values for data rows 2..11
$saveValues = @(0, 0, 1, 0, 0, 0, 0, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
